$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.616.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.674.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.45%  "

$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3948"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3948"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.88%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.399"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08646"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.22%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.316"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001319"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.698"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.37%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.679.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.55%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06999"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.089"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.00%  "

$ws.Range("E22").Value = "  -0.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.619.79"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.344"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.773"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.839"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "160.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "145.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.288"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.499"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.866.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.03079"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08274"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.930"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2813"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9940"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09620"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.515"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.51%  "

$ws.Range("E41").Value = "  -5.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7905"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.564"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7106"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.173"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08666"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.95%  "

$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("E50").Value = "  -4.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.17%  "
